$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value
$ws.Range("B3").Value = "6.0.0"

# Update Date value
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Set Publisher value (was empty)
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Old row 11 (duplicate "Contact" row) is removed entirely, shifting rows 12-15 up to 11-14
$ws.Rows.Item(11).Delete()
